$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function TitleCaseConnectors($s) {
    if ($s.Equals("TOTAL")) {
        return "Total"
    }
    $words = $s.Split(" ")
    $out = @()
    foreach ($w in $words) {
        if ($w.Equals("de")) { $out += "De" }
        elseif ($w.Equals("del")) { $out += "Del" }
        elseif ($w.Equals("la")) { $out += "La" }
        elseif ($w.Equals("las")) { $out += "Las" }
        elseif ($w.Equals("los")) { $out += "Los" }
        elseif ($w.Equals("el")) { $out += "El" }
        else { $out += $w }
    }
    return [string]::Join(" ", $out)
}

# Rename header columns to the new (English/snake_case) names
$ws.Range("A1").Value() = "mx_state"
$ws.Range("B1").Value() = "mx_municipality"
$ws.Range("C1").Value() = "n_matriculas"
$ws.Range("D1").Value() = "pct_matriculas"

# Title-case the connector words (de/del/la/las/los/el) in the state (A)
# and municipality (B) columns, and normalize "TOTAL" -> "Total".
for ($r = 2; $r -le 333; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -ne $null) {
            $newv = TitleCaseConnectors $v
            if (-not $newv.Equals($v)) {
                $cell.Value() = $newv
            }
        }
    }
}

# Remove the trailing footer/metadata rows (335-339)
$ws.Rows("335:339").Delete()
